$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D price cells to remain text (they contain dot-grouped
# numbers like "26.083.30" which Excel would otherwise coerce to a number).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.083.30"
$ws.Range("E2").Value = "  +1.14%  "

$ws.Range("D3").Value = "1.767.81"
$ws.Range("E3").Value = "  +1.31%  "

$ws.Range("D5").Value = "237.73"
$ws.Range("E5").Value = "  -0.36%  "

$ws.Range("E6").Value = "  +0.19%  "

$ws.Range("D7").Value = "0.5226"
$ws.Range("E7").Value = "  +3.55%  "

$ws.Range("D8").Value = "0.2755"
$ws.Range("E8").Value = "  +0.93%  "

$ws.Range("D9").Value = "40.46"
$ws.Range("E9").Value = "  -3.67%  "

$ws.Range("D10").Value = "0.06193"
$ws.Range("E10").Value = "  +0.64%  "

$ws.Range("D11").Value = "1.776.22"
$ws.Range("E11").Value = "  +1.72%  "

$ws.Range("D12").Value = "0.07019"
$ws.Range("E12").Value = "  +1.34%  "

$ws.Range("D13").Value = "15.70"
$ws.Range("E13").Value = "  +1.21%  "

$ws.Range("D14").Value = "0.6438"
$ws.Range("E14").Value = "  +6.92%  "

$ws.Range("D15").Value = "4.532"

$ws.Range("D16").Value = "78.09"

$ws.Range("D17").Value = "1.002"
$ws.Range("E17").Value = "  +0.21%  "

$ws.Range("E18").Value = "  +0.19%  "

$ws.Range("D19").Value = "26.093.13"
$ws.Range("E19").Value = "  +1.17%  "

$ws.Range("D20").Value = "11.66"
$ws.Range("E20").Value = "  +0.27%  "

$ws.Range("D21").Value = "0.000006746"
$ws.Range("E21").Value = "  -2.05%  "

$ws.Range("D22").Value = "2.000.72"
$ws.Range("E22").Value = "  +1.64%  "

$ws.Range("D23").Value = "4.073"
$ws.Range("E23").Value = "  +0.33%  "

$ws.Range("D24").Value = "8.448"
$ws.Range("E24").Value = "  +3.36%  "

$ws.Range("D25").Value = "5.186"
$ws.Range("E25").Value = "  -1.06%  "

$ws.Range("D26").Value = "138.88"
$ws.Range("E26").Value = "  +0.68%  "

$ws.Range("D27").Value = "1.489"
$ws.Range("E27").Value = "  +1.57%  "

$ws.Range("D28").Value = "1.851"
$ws.Range("E28").Value = "  +1.80%  "

$ws.Range("D29").Value = "15.21"
$ws.Range("E29").Value = "  +1.33%  "

$ws.Range("D30").Value = "103.25"
$ws.Range("E30").Value = "  -0.58%  "

$ws.Range("D31").Value = "0.08402"
$ws.Range("E31").Value = "  +3.26%  "

$ws.Range("D32").Value = "3.700"
$ws.Range("E32").Value = "  -0.22%  "

$ws.Range("D33").Value = "3.451"
$ws.Range("E33").Value = "  -0.98%  "

$ws.Range("D34").Value = "0.04454"
$ws.Range("E34").Value = "  -2.00%  "

$ws.Range("D35").Value = "2.614"
$ws.Range("E35").Value = "  -0.01%  "

$ws.Range("E36").Value = "  +1.65%  "

$ws.Range("D37").Value = "0.6037"
$ws.Range("E37").Value = "  -1.11%  "

$ws.Range("E38").Value = "  +2.72%  "

$ws.Range("D39").Value = "0.01590"
$ws.Range("E39").Value = "  +2.26%  "

$ws.Range("D40").Value = "1.984"
$ws.Range("E40").Value = "  +2.41%  "

$ws.Range("D41").Value = "1.002"
$ws.Range("E41").Value = "  +0.30%  "

$ws.Range("D42").Value = "102.64"
$ws.Range("E42").Value = "  +0.70%  "

$ws.Range("E43").Value = "  +0.66%  "

$ws.Range("D44").Value = "0.7508"
$ws.Range("E44").Value = "  +1.71%  "

$ws.Range("D45").Value = "4.936"
$ws.Range("E45").Value = "  -0.67%  "

$ws.Range("D46").Value = "0.05518"
$ws.Range("E46").Value = "  +2.65%  "

$ws.Range("D47").Value = "6.338"
$ws.Range("E47").Value = "  +6.10%  "

$ws.Range("E48").Value = "  +0.34%  "

$ws.Range("D49").Value = "30.19"
$ws.Range("E49").Value = "  -0.02%  "

$ws.Range("D50").Value = "52.54"
$ws.Range("E50").Value = "  -0.07%  "

$ws.Range("D51").Value = "1.003"
$ws.Range("E51").Value = "  +0.80%  "

# Restore the default (unformatted) style on column D so only the cell
# contents change, matching the original styling (no explicit style index).
$ws.Range("D2:D51").Style = "Normal"
